$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Function/Description columns for the new "Driver" class rows
$ws.Range("B27").Value = "readFile(String filename) "
$ws.Range("C27").Value = "reads strings from given file (filename) "

$ws.Range("B28").Value = "getUsers(ArrayList<String> userStrings) "
$ws.Range("B29").Value = "getItems(ArrayList<String>itemStrings) "
$ws.Range("C29").Value = "takes in an ArrayList of Strings and builds an ArrayList<Items> "
$ws.Range("C28").Value = "takes in an ArrayList of Strings and builds an ArrayList<Users> "

$ws.Range("B30").Value = "writeToFile(ArrayList<String> strings, String filename) "
$ws.Range("C30").Value = "writes Strings from ArrayList: strings to file: filename"

# Fill in the Class column for the new rows
$ws.Range("A27").Value = "Driver"
$ws.Range("A28").Value = "Driver"
$ws.Range("A29").Value = "Driver"
$ws.Range("A30").Value = "Driver"

# Widen columns A and B to fit the new, longer content
$ws.Columns.Item(1).ColumnWidth = 26.666666666666664
$ws.Columns.Item(2).ColumnWidth = 49.666666666666664

# Update the active selection
$ws.Range("B33").Select() | Out-Null
